$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "58.412.22"
$ws.Range("E2").Value = "  -0.38%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.141.82"
$ws.Range("E3").Value = "  +1.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D5").Value = "'532.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "

# Row 6 - Solana (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D6").Value = "'143.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.141.07"
$ws.Range("E8").Value = "  +1.22%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.80%  "

# Row 10 - Toncoin (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D10").Value = "'7.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.76%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.40%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +2.52%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.681.22"
$ws.Range("E13").Value = "  +1.15%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +3.22%  "

# Row 15 - Avalanche (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D15").Value = "'25.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.53%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -0.20%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "58.424.57"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.135.18"
$ws.Range("E18").Value = "  +1.19%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -0.44%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.50%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.34%  "

# Row 22 - BitcoinCash (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D22").Value = "'344.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.88%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.11%  "

# Row 24 - Polygon (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D24").Value = "'0.513"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.43%  "

# Row 25 - Litecoin (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D25").Value = "'67.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.58%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.66%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.43%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0935"
$ws.Range("E28").Value = "  +1.90%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("E29").Value = "  +3.37%  "

# Row 31 - RenderToken (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D31").Value = "'6.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.82%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +1.52%  "

# Row 33 - EthereumClassic (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D33").Value = "'21.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "

# Row 34 - Fetch.AI
$ws.Range("E34").Value = "  -0.77%  "

# Row 35 - was Monero, now NEARProtocol (numeric-looking text: quote-prefix then restore default style)
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.30%  "

# Row 36 - was NEARProtocol, now Monero (numeric-looking text: quote-prefix then restore default style)
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'158.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.70%  "

# Row 37 - Aptos
$ws.Range("E37").Value = "  +2.71%  "

# Row 38 - EnergySwap (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D38").Value = "'26.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.51%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -4.13%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +10.42%  "

# Row 41 - Hedera
$ws.Range("E41").Value = "  -1.07%  "

# Row 42 - Mantle (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D42").Value = "'0.709"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.71%  "

# Row 43 - Filecoin (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D43").Value = "'4.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.68%  "

# Row 44 - RenzoRestakedETH
$ws.Range("D44").Value = "3.181.04"
$ws.Range("E44").Value = "  +1.08%  "

# Row 45 - OKB (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D45").Value = "'36.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.02%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  -0.04%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +2.52%  "

# Row 48 - Maker
$ws.Range("D48").Value = "2.287.62"
$ws.Range("E48").Value = "  -0.52%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  +4.32%  "

# Row 50 - InjectiveProtocol (numeric-looking text: quote-prefix then restore default style)
$ws.Range("D50").Value = "'20.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.38%  "

# Row 51 - Cosmos
$ws.Range("E51").Value = "  +1.66%  "
